$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F4").Value = 1464
    $ws.Range("G6").Value = 65
    $ws.Range("F7").Value = 110
    $ws.Range("F8").Value = 42
    $ws.Range("F9").Value = 254
}
